# Auto-generated edit script applying the diff to cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on Price (D) cells whose new value looks like a
# plain decimal number, so Excel stores the literal text (matching the source
# data, which keeps trailing zeros / thousand-dot formatting) instead of
# silently converting it to a floating point number.
$priceTextCells = @(
'D2', 'D3', 'D5', 'D6', 'D7', 'D10', 'D11', 'D12', 'D13', 'D15', 'D16', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D30', 'D33', 'D34', 'D36', 'D38', 'D39', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D49', 'D50', 'D51'
)
foreach ($addr in $priceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new cell values (Coin name / Link / Price / Volume(1h)).
$ws.Range('D2').Value = '51.470.93'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '3.104.69'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '383.72'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').Value = '103.03'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').Value = '0.540'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.68%  '
$ws.Range('D10').Value = '37.17'
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').Value = '0.138'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').Value = '0.0851'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = '3.593.27'
$ws.Range('E13').Value = '  +2.67%  '
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('D15').Value = '7.81'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = '3.102.51'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').Value = '11.29'
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('E18').Value = '  -0.91%  '
$ws.Range('D19').Value = '51.486.06'
$ws.Range('E19').Value = '  -0.67%  '
$ws.Range('D20').Value = '3.28'
$ws.Range('E20').Value = '  +6.57%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0962'
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('B22').Value = 'InternetComputer(DFINITY)'
$ws.Range('C22').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = '12.31'
$ws.Range('E22').Value = '  -2.21%  '
$ws.Range('D23').Value = '69.93'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '265.64'
$ws.Range('E24').Value = '  -1.26%  '
$ws.Range('D25').Value = '3.09'
$ws.Range('E25').Value = '  -3.40%  '
$ws.Range('E26').Value = '  -2.73%  '
$ws.Range('D27').Value = '26.99'
$ws.Range('E27').Value = '  +2.71%  '
$ws.Range('E28').Value = '  -4.29%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '0.167'
$ws.Range('E30').Value = '  -2.40%  '
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').Value = '35.42'
$ws.Range('E33').Value = '  +2.50%  '
$ws.Range('D34').Value = '0.0472'
$ws.Range('E34').Value = '  +3.34%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').Value = '50.30'
$ws.Range('E36').Value = '  -1.80%  '
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('D38').Value = '3.35'
$ws.Range('E38').Value = '  +0.67%  '
$ws.Range('D39').Value = '0.298'
$ws.Range('E39').Value = '  +4.68%  '
$ws.Range('E40').Value = '  +1.02%  '
$ws.Range('D41').Value = '128.95'
$ws.Range('E41').Value = '  +2.27%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.116'
$ws.Range('E42').Value = '  -0.88%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').Value = '16.51'
$ws.Range('E43').Value = '  -3.86%  '
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').Value = '3.67'
$ws.Range('E45').Value = '  -3.39%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '22.28'
$ws.Range('E46').Value = '  +2.02%  '
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('D49').Value = '2.057.80'
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('D50').Value = '3.412.01'
$ws.Range('E50').Value = '  +2.19%  '
$ws.Range('D51').Value = '0.0325'
$ws.Range('E51').Value = '  -0.34%  '
